$wb = $excel.ActiveWorkbook

# Device sheet: update uptime string
$wsDevice = $wb.Worksheets.Item("Device")
$wsDevice.Range("G2").Value = "20 hours, 34 minutes"

# Mem_CPU sheet: update memory usage figures
$wsMem = $wb.Worksheets.Item("Mem_CPU")
$wsMem.Range("C3").Value = 405610408
$wsMem.Range("D3").Value = 1282750376

# Buffer sheet: update buffer count figures
$wsBuffer = $wb.Worksheets.Item("Buffer")
$wsBuffer.Range("C2").Value = 314690
$wsBuffer.Range("C3").Value = 304127
$wsBuffer.Range("C4").Value = 68852
$wsBuffer.Range("C5").Value = 44516
